# Updated time sheets and some minor changes to report
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Week 12 (col Q) updates: new/changed entries get the elapsed-time
#     format "[h]:mm" (matches the other populated weekly-hours cells) ---
$ws.Range("Q2").Value = 1.2291666666666667
$ws.Range("Q2").NumberFormat = "[h]:mm"

$ws.Range("Q3").Value = 1.6875

$ws.Range("Q4").Value = 1.5
$ws.Range("Q4").NumberFormat = "[h]:mm"

$ws.Range("Q5").Value = 1.875
$ws.Range("Q5").NumberFormat = "[h]:mm"

$ws.Range("Q6").Value = 0.041666666666666664
$ws.Range("Q6").NumberFormat = "[h]:mm"

$ws.Range("Q7").Value = 1.3333333333333333
$ws.Range("Q7").NumberFormat = "[h]:mm"

$ws.Range("Q8").Value = 1.4583333333333333
$ws.Range("Q8").NumberFormat = "[h]:mm"

# --- Week 13 (col R) new entry for row 3 ---
$ws.Range("R3").Value = 0.625

# --- Row 9 (Tim Anderson) newly-logged weekly hours ---
$ws.Range("F9").Value = 0.041666666666666664
$ws.Range("G9").Value = 0.041666666666666664
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("K9").Value = 0.08333333333333333
$ws.Range("L9").Value = 0.1875
$ws.Range("M9").Value = 0.08333333333333333

$ws.Range("Q9").Value = 1.1041666666666667
$ws.Range("Q9").NumberFormat = "[h]:mm"

$ws.Range("R9").Value = 0.16666666666666666

# --- Starting hours bumped ---
$ws.Range("B11").Value = 3.3333333333333335

# Recalculate so dependent formulas (C2:C9, D2:D9, D11) pick up new values
$excel.Calculate()

# Restore the active selection as left by the author
$ws.Range("M9").Select() | Out-Null
